# Append a new data row (row 26) to the CityResaleNum sheet, mirroring the
# structure of the existing rows: columns A-D hold text (date, time,
# weekday, week-number-as-text) and columns E-T hold numeric resale counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 26
$prev = 2   # an existing data row used as a style/format template

# Columns A, B and D look like a date / time / number but must stay plain
# text (matching every other row in the sheet), so force a text format
# before writing the values to stop Excel auto-converting them.
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("B$row").NumberFormat = "@"
$ws.Range("D$row").NumberFormat = "@"

$ws.Range("A$row").Value = "2023-06-06"
$ws.Range("B$row").Value = "21:47:27"
$ws.Range("C$row").Value = "Tuesday"
$ws.Range("D$row").Value = "23"

$ws.Range("E$row").Value = 120538
$ws.Range("F$row").Value = 134212
$ws.Range("G$row").Value = 159769
$ws.Range("H$row").Value = 130501
$ws.Range("I$row").Value = 175141
$ws.Range("J$row").Value = 112572
$ws.Range("K$row").Value = 200471
$ws.Range("L$row").Value = 220290
$ws.Range("M$row").Value = 172380
$ws.Range("N$row").Value = 119594
$ws.Range("O$row").Value = 38452
$ws.Range("P$row").Value = 34615
$ws.Range("Q$row").Value = 50488
$ws.Range("R$row").Value = -1
$ws.Range("S$row").Value = 36846
$ws.Range("T$row").Value = -1

# Re-apply the same (default) cell style the rest of the table uses so the
# new row doesn't end up with a stray "Text" number format applied to it.
$ws.Range("A$row").Style = $ws.Range("A$prev").Style
$ws.Range("B$row").Style = $ws.Range("B$prev").Style
$ws.Range("C$row").Style = $ws.Range("C$prev").Style
$ws.Range("D$row").Style = $ws.Range("D$prev").Style
